$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hostname value (B7): adeye03u -> adeye06u
$ws.Range("B7").Value = "adeye06u"

# Update ipaddress value (B2): 130.237.10.123 -> 192.168.122.1
$ws.Range("B2").Value = "192.168.122.1"

# Move selection to B2
$ws.Range("B2").Select()
